# Update Name of Algo
# Apply updated numeric values (result data) to Sheet1 as described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value  = 7.366
$ws.Range("A9").Value  = -21.462
$ws.Range("B9").Value  = 6.456
$ws.Range("B11").Value = 6.414
$ws.Range("A18").Value = -21.995
$ws.Range("A20").Value = -20.623
$ws.Range("B23").Value = 7.731
$ws.Range("B24").Value = 5.571000000000001
$ws.Range("B26").Value = 5.752
$ws.Range("A27").Value = -21.637
$ws.Range("B34").Value = 7.175
$ws.Range("A35").Value = -21.842
$ws.Range("B35").Value = 5.794000000000001
$ws.Range("B48").Value = 5.616
$ws.Range("B49").Value = 6.237
$ws.Range("B52").Value = 5.573
$ws.Range("B66").Value = 5.202
$ws.Range("B67").Value = 5.383
$ws.Range("A69").Value = -21.291
$ws.Range("A76").Value = -20.392
$ws.Range("A78").Value = -20.743
$ws.Range("B78").Value = 6.742
$ws.Range("B80").Value = 8.298999999999999
$ws.Range("A82").Value = -21.81
$ws.Range("A83").Value = -21.509
$ws.Range("A93").Value = -21.533
$ws.Range("B99").Value = 5.3
$ws.Range("B104").Value = 7.255
